$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 342, shifting existing rows 342..443 down to 343..444
$ws.Rows.Item(342).Insert()

# Populate the newly inserted row 342 with the new weekly price record
$ws.Range("A342").Value = 11
$ws.Range("B342").Value = "Vega Monumental Concepción"
$ws.Range("C342").Value = "Bíobío"
$ws.Range("D342").Value = 44988
$ws.Range("E342").Value = 8
$ws.Range("F342").Value = 100112023
$ws.Range("G342").Value = "Brócoli"
$ws.Range("H342").Value = "Sin especificar"
$ws.Range("I342").Value = "Primera"
$ws.Range("J342").Value = 2500
$ws.Range("K342").Value = 900
$ws.Range("L342").Value = 1000
$ws.Range("M342").Value = 960
$ws.Range("N342").Value = "$/unidad"
$ws.Range("O342").Value = "Región Metropolitana"
$ws.Range("P342").Value = 960
$ws.Range("Q342").Value = 1
$ws.Range("R342").Value = "Hortaliza"
